$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.604.35"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "2.292.19"
$ws.Range("E3").Value = "  -1.23%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "495.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "2.290.08"
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0949"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("E13").Value = "  -2.92%  "
$ws.Range("D14").Value = "2.696.54"
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("D16").Value = "54.510.66"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "2.292.44"
$ws.Range("E18").Value = "  -1.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.44%  "
$ws.Range("E20").Value = "  +2.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "305.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.76%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  -1.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.02%  "
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.151"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.62%  "
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("D29").Value = "2.392.39"
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.59%  "
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("D33").Value = "0.0₃0686"
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("E34").Value = "  +2.52%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("E39").Value = "  +3.13%  "
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "35.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("E44").Value = "  +2.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "128.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.51%  "
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0893"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.549"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "242.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0482"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.60%  "
